# Update power-line flow values for the "380 kV" case (rows 2-25) in Sheet1.
# Columns updated: B, D, E, F, G, H, I, M (columns C, J, K, L, N, O stay 0).
# Column order in each data row below: Row, B, D, E, F, G, H, I, M

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.4271489846939573, 0.004614702950672012, 0.135093437805935, 0.8853024841834696, 0.8394016251781977, 0.6505700499109253, 0.9272109797212238, 1.801843627829953),
    @(3, 0.3754131897714501, 0.004035205983424106, 0.1282488325481737, 0.8148567801142974, 0.7563778988490242, 0.6197482330507, 0.8491968225310558, 1.577287430541674),
    @(4, 0.3435217106739969, 0.003694052512262402, 0.1241556180976957, 0.7725369654360605, 0.7062663125626045, 0.6015454560000251, 0.8016903550215631, 1.439033166115138),
    @(5, 0.3304945580995593, 0.003558522972419098, 0.1225155263798854, 0.7555202522952271, 0.6860558564297889, 0.5943059865521434, 0.7824311305090532, 1.382599662917983),
    @(6, 0.3283295436708897, 0.003536223772130143, 0.1222448894276553, 0.7527082959238101, 0.682712426399803, 0.5931145539471743, 0.7792392201439071, 1.373223284273593),
    @(7, 0.3433461473186696, 0.003692210822467956, 0.1241333856168012, 0.7723065528628155, 0.7059929047117066, 0.6014471042494733, 0.8014302121140418, 1.43827246273392),
    @(8, 0.409336766479413, 0.004411732776475219, 0.132710957134087, 0.8608153106412288, 0.8105910874242284, 0.6397908836927684, 0.9002302557032777, 1.724494993633186),
    @(9, 0.5377370291978423, 0.005947846494301245, 0.1503826086062361, 1.042071531404417, 1.022909184827881, 0.7208624866065065, 1.097086580383547, 2.282796493050995),
    @(10, 0.6314550549882938, 0.007165071012575908, 0.1638639064313168, 1.180331179345984, 1.183771661865791, 0.7842278690313265, 1.243608395439878, 2.691218629851051),
    @(11, 0.6739562360519074, 0.007740785194084054, 0.1701008075289252, 1.244423546862834, 1.258116802150539, 0.8139282583092324, 1.310678010012339, 2.876655317755706),
    @(12, 0.6900312825355286, 0.007962181435544835, 0.1724771677637804, 1.268872970717553, 1.286445885769297, 0.8253047541907677, 1.336135342880539, 2.94682481869566),
    @(13, 0.6865700932559093, 0.007914345614544516, 0.1719647338828167, 1.263599276742639, 1.280336741899845, 0.8228487975814005, 1.330650007133926, 2.931714844249569),
    @(14, 0.6752791267635416, 0.007758930447774759, 0.1702960221630576, 1.246431387845774, 1.260443874071825, 0.8148615902798326, 1.31277120852215, 2.882429233728345),
    @(15, 0.6683605713373311, 0.007664181790634217, 0.1692757751510783, 1.235939081385936, 1.248282110429813, 0.8099861859504074, 1.301827672743258, 2.85223368790929),
    @(16, 0.6286748277079823, 0.007127908280875062, 0.1634583723955458, 1.176167269745463, 1.178937236044248, 0.7823048087493589, 1.239233582534922, 2.679092720778215),
    @(17, 0.6042949722744027, 0.006804712338809082, 0.159915982253473, 1.139810713129293, 1.13670136969003, 0.765549981230464, 1.200940455419044, 2.572784910497035),
    @(18, 0.5902599513465816, 0.00662087504933595, 0.1578883236587245, 1.119011667472506, 1.112517852641474, 0.7559951930807358, 1.178954439072086, 2.511605899439274),
    @(19, 0.5855058191599483, 0.006558977363230412, 0.1572034940399405, 1.111988545908105, 1.104348276357342, 0.7527741011766125, 1.171517091709902, 2.490885976699388),
    @(20, 0.606891532354723, 0.006838902656355117, 0.1602920612399856, 1.14366926038042, 1.141186064312762, 0.7673250324821765, 1.205012772033328, 2.584105033825466),
    @(21, 0.6785960813943177, 0.007804485962438434, 0.1707857706810643, 1.251469098770457, 1.266282043418272, 0.8172040790163635, 1.318021031495192, 2.896907005549281),
    @(22, 0.7253469575566101, 0.008455391768777076, 0.1777288251507869, 1.322969347187552, 1.349070404676496, 0.8505601367290865, 1.392225354192107, 3.101041966647244),
    @(23, 0.7004054908222201, 0.008106100825841622, 0.1740155571284987, 1.284710238004863, 1.304787612220252, 0.8326868367687155, 1.352589444448142, 2.992118653519555),
    @(24, 0.605717685818945, 0.006823439087337846, 0.160122008151049, 1.14192449141467, 1.139158227860804, 0.7665222898106947, 1.203171586676291, 2.578987396792883),
    @(25, 0.5031096965234099, 0.005517625753427069, 0.145513032004061, 0.9921723771782638, 0.9646511563996683, 0.6982787314453276, 1.043499664472861, 2.132075121552504)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]   # B
    $ws.Cells.Item($r, 4).Value  = $row[2]   # D
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G
    $ws.Cells.Item($r, 8).Value  = $row[6]   # H
    $ws.Cells.Item($r, 9).Value  = $row[7]   # I
    $ws.Cells.Item($r, 13).Value = $row[8]   # M
}
